$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Establish new shared strings in the same order as the target workbook:
#  24 -> "Negative Cs, Counts below detection" (capital C)
#  25 -> "Negative Cs, counts below detection" (lowercase c)
#  26 -> "Negative Cs, setting to 0"

# Row 3: Comment -> "Negative Cs, Counts below detection"
$ws.Range("N3").Value = "Negative Cs, Counts below detection"

# Row 2: Comment -> "Negative Cs, counts below detection"
$ws.Range("N2").Value = "Negative Cs, counts below detection"

# Row 4: clean-up of sorption computations, later date, included
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("L4").Value = Get-Date -Year 2016 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("M4").Value = $true
$ws.Range("N4").Value = "Negative Cs, setting to 0"

# Row 7 & 8: Comment -> "Negative Cs, counts below detection"
$ws.Range("N7").Value = "Negative Cs, counts below detection"
$ws.Range("N8").Value = "Negative Cs, counts below detection"

# Update the active selection in the frozen (bottom-right) pane to reflect
# the new scroll/selection position, without disturbing the existing freeze
# (xSplit=2 / ySplit=1 stays as-is).
$ws.Range("N9").Select()
